# Applies updated odds/statistics values for rows 8, 9 and 10 on Sheet1
# as described by the commit "Atualizando o arquivo XLSX".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 8 ---
$ws.Range("G8").Value  = 2.55
$ws.Range("H8").Value  = 3.25
$ws.Range("I8").Value  = 2.7
$ws.Range("M8").Value  = 1.06
$ws.Range("N8").Value  = 8
$ws.Range("O8").Value  = 1.3
$ws.Range("P8").Value  = 3.4
$ws.Range("Q8").Value  = 2
$ws.Range("R8").Value  = 1.8
$ws.Range("U8").Value  = 1.8
$ws.Range("V8").Value  = 1.91
$ws.Range("AB8").Value = 29
$ws.Range("AC8").Value = 9.5
$ws.Range("AD8").Value = 6.5
$ws.Range("BA8").Value = 67

# --- Row 9 ---
$ws.Range("G9").Value  = 1.55
$ws.Range("I9").Value  = 6
$ws.Range("J9").Value  = 2.05
$ws.Range("U9").Value  = 1.67
$ws.Range("V9").Value  = 2.1
$ws.Range("W9").Value  = 9
$ws.Range("X9").Value  = 8.5
$ws.Range("AE9").Value = 15
$ws.Range("AG9").Value = 151
$ws.Range("AH9").Value = 19
$ws.Range("AJ9").Value = 19
$ws.Range("AP9").Value = 17
$ws.Range("AS9").Value = 101
$ws.Range("AU9").Value = 8
$ws.Range("AV9").Value = 51
$ws.Range("AY9").Value = 29
$ws.Range("BA9").Value = 101

# --- Row 10 ---
$ws.Range("K10").Value  = 2.25
$ws.Range("M10").Value  = 1.04
$ws.Range("N10").Value  = 13
$ws.Range("O10").Value  = 1.22
$ws.Range("P10").Value  = 4
$ws.Range("Q10").Value  = 1.75
$ws.Range("R10").Value  = 2.05
$ws.Range("U10").Value  = 1.62
$ws.Range("V10").Value  = 2.2
$ws.Range("W10").Value  = 9
$ws.Range("AC10").Value = 12
$ws.Range("AH10").Value = 12
$ws.Range("AM10").Value = 29
